# Add a "Status" column (F) to the API sheet, with a couple of "Done"
# markers, as described in the commit: "Added Status column to API so
# to see what's ready."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Done" markers in column F (added first so "Done" gets the lower
#     shared-string index, matching the original author's edit order) ---
$ws.Range("F2").Value = "Done"
$ws.Range("F2").Font.Color = 0
$ws.Range("F2").WrapText = $true
$ws.Range("F2").VerticalAlignment = -4108

$ws.Range("F10").Value = "Done"
$ws.Range("F10").Font.Color = 0
$ws.Range("F10").WrapText = $true
$ws.Range("F10").VerticalAlignment = -4108

$ws.Range("F11").Value = "Done"

# --- Header cell F1: "Status" -----------------------------------------
$ws.Range("F1").Value = "Status"
$ws.Range("F1").Font.Bold = $true

# --- Keep the previously-selected cell in the same relative spot -------
$ws.Range("E5").Select() | Out-Null
